$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.78"
$ws.Range("E2").Value = "'1.24%"

$ws.Range("D3").Value = "'27.33"
$ws.Range("E3").Value = "'1.95%"

$ws.Range("D4").Value = "'4.885"
$ws.Range("E4").Value = "'-0.73%"

$ws.Range("D5").Value = "'0.06364"
$ws.Range("E5").Value = "'0.35%"

$ws.Range("E6").Value = "'0.25%"

$ws.Range("D7").Value = "'1.275"
$ws.Range("E7").Value = "'-10.82%"

$ws.Range("D8").Value = "'0.8819"
$ws.Range("E8").Value = "'-0.62%"

$ws.Range("E9").Value = "'3.22%"

$ws.Range("E10").Value = "'2.08%"

$ws.Range("D11").Value = "'0.07529"
$ws.Range("E11").Value = "'2.02%"

$ws.Range("D12").Value = "'0.02979"
$ws.Range("E12").Value = "'-5.96%"

$ws.Range("D13").Value = "'0.09010"
$ws.Range("E13").Value = "'-0.58%"

$ws.Range("D14").Value = "'0.001575"
$ws.Range("E14").Value = "'0.38%"

$ws.Range("D15").Value = "'0.0006387"
$ws.Range("E15").Value = "'0.75%"

$ws.Range("D16").Value = "'0.006032"
$ws.Range("E16").Value = "'-0.04%"

$ws.Range("D17").Value = "'3.459"
$ws.Range("E17").Value = "'-0.71%"

$ws.Range("E18").Value = "'-1.38%"

$ws.Range("D19").Value = "'2.285"
$ws.Range("E19").Value = "'0.18%"

$ws.Range("D21").Value = "'0.1337"
$ws.Range("E21").Value = "'0.11%"

$ws.Range("D22").Value = "'3.913"
$ws.Range("E22").Value = "'0.11%"

$ws.Range("D23").Value = "'0.04415"
$ws.Range("E23").Value = "'1.66%"

$ws.Range("D24").Value = "'0.001172"
$ws.Range("E24").Value = "'-0.41%"

$ws.Range("E25").Value = "'5.96%"

$ws.Range("E26").Value = "'-0.22%"

$ws.Range("E27").Value = "'-0.44%"

$ws.Range("D40").Value = "'0.04146"
$ws.Range("E40").Value = "'2.88%"

$ws.Range("D41").Value = "'0.006843"
$ws.Range("E41").Value = "'3.21%"

$ws.Range("D42").Value = "'0.1179"
$ws.Range("E42").Value = "'1.02%"

$ws.Range("D43").Value = "'0.002020"
$ws.Range("E43").Value = "'-14.59%"

$ws.Range("E44").Value = "'-10.91%"

$ws.Range("D45").Value = "'0.00005171"
$ws.Range("E45").Value = "'-1.96%"

$ws.Range("D47").Value = "'0.02024"
$ws.Range("E47").Value = "'-4.91%"
